$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value as TEXT (matching the source data,
# which stores these figures as shared-string text, not numbers) while
# keeping the cell's original (default/general) style - an apostrophe
# prefix forces text entry, then resetting the style back to "Normal"
# drops the transient quote-prefix style Excel would otherwise apply.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 13 - Enterprises density (per 1000 people): Micro / SMEs / MSMEs
Set-TextValue $ws.Range("B13") "85.39"
Set-TextValue $ws.Range("C13") "3.91"
Set-TextValue $ws.Range("D13") "89.31"

# Row 14 - Employment (% of total): Micro / SMEs / MSMEs
Set-TextValue $ws.Range("B14") "30.11"
Set-TextValue $ws.Range("C14") "38.36"
Set-TextValue $ws.Range("D14") "68.47"

# Row 16 - Enterprises (% of total): Micro / SMEs / MSMEs
Set-TextValue $ws.Range("B16") "95.47"
Set-TextValue $ws.Range("C16") "4.38"
Set-TextValue $ws.Range("D16") "99.84"
